# fall 24 week 3 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D2").Value = 10.96

$ws.Range("D3").Value = 10.09
$ws.Range("F3").Value = 10.43

$ws.Range("B4").Value = 9.039999999999999
$ws.Range("C4").Value = 9.91
$ws.Range("E4").Value = 10.63
$ws.Range("F4").Value = 9.720000000000001
$ws.Range("J4").Value = 11.5

$ws.Range("D5").Value = 9.369999999999999
$ws.Range("F5").Value = 10.18

$ws.Range("C6").Value = 9.57
$ws.Range("D6").Value = 10.28
$ws.Range("E6").Value = 9.82
$ws.Range("G6").Value = 10.51

$ws.Range("F7").Value = 9.49
$ws.Range("H7").Value = 10.02

$ws.Range("G8").Value = 9.98

$ws.Range("D10").Value = 8.5
